# Refresh the cryptocurrency price / 1h-volume table with the latest
# scrape, including the WrappedEther/WrappedBTC rows swapping ranking
# positions (rows 17 and 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.802.67"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "3.661.62"
$ws.Range("E3").Value = "  +7.75%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "3.655.30"
$ws.Range("E7").Value = "  +7.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.611"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "682.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "4.245.32"
$ws.Range("E15").Value = "  +7.45%  "
$ws.Range("E16").Value = "  +4.26%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "71.937.79"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.651.03"
$ws.Range("E18").Value = "  +7.28%  "
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.943"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.50%  "
$ws.Range("E24").Value = "  +2.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  +2.71%  "
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "578.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "3.737.64"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "0.0₃0763"
$ws.Range("E41").Value = "  +3.41%  "
$ws.Range("E42").Value = "  +4.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0468"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("E47").Value = "  +3.17%  "
$ws.Range("E48").Value = "  +3.98%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.30%  "
